# Generate Report for Handoff
#
# Refresh the localization-status report with a new handoff run: two
# tracked source files are replaced by a new pair, status/timestamps are
# updated, and the redundant "Source File Name" / "Latest Handoff File"
# duplicate columns (F/G) are dropped from the per-locale sheets.

$wb = $excel.ActiveWorkbook

$newName1 = "b496a5f4-4a97-4f66-b179-2798b54a02f9.md"
$newName2 = "ffff93b24ac0-24d3-45af-9f08-77195abee18a.md"
$newStatus = "Ready for handoff"
$newHandoffDate = "2016-51-17 20:51:03"

$xlfBase = "b496a5f4-4a97-4f66-b179-2798b54a02f9.f4805f986f0de95a7dd32215cb1f767a935043db"
$newHandoffDatetime = "2016-03-17 20:50:56"
$newHandbackDatetime = "0001-01-01 00:00:00"
$newDeHandoffDatetime = "2016-03-17 20:51:03"

# Hyperlink target URLs are untouched by this edit (the underlying commit
# references do not change) -- only the visible cell text is refreshed.
$md1Url = "https://github.com/OpenLocalizationTest/oltest/blob/77105346726004afb18d8ecda6fb6c9e25121a62/e2e/3edd96af-b0cd-414c-ba2b-54cb1bb2ebd7.md"
$md2Url = "https://github.com/OpenLocalizationTest/oltest/blob/77105346726004afb18d8ecda6fb6c9e25121a62/e2e/58c7fef4-7043-4af8-9658-444a7f9a32d6.md"
$zhXlf1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/213f47510437eeb4cc44cc31d3731c96f9cb08f0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3edd96af-b0cd-414c-ba2b-54cb1bb2ebd7.ed74af1e17c37847078fbf243195a30a412ec1b6.zh-cn.xlf"
$zhXlf2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/213f47510437eeb4cc44cc31d3731c96f9cb08f0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/58c7fef4-7043-4af8-9658-444a7f9a32d6.5e2c6408b1e4467cfb04aeec1188a48f0bf3abf9.zh-cn.xlf"
$deXlf1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9bf686bf731d2358bd4f0ad9ae5bf5076d8c9b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3edd96af-b0cd-414c-ba2b-54cb1bb2ebd7.ed74af1e17c37847078fbf243195a30a412ec1b6.de-de.xlf"
$deXlf2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9bf686bf731d2358bd4f0ad9ae5bf5076d8c9b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/58c7fef4-7043-4af8-9658-444a7f9a32d6.5e2c6408b1e4467cfb04aeec1188a48f0bf3abf9.de-de.xlf"

# BGR-ordered value of RGB FF6495ED -- the workbook's "HyperLink" font color.
$hyperlinkColor = 15570276

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $hyperlinkColor
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
}

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newName1
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("D2").Value = $newHandoffDate

$wsOverview.Range("A3").Value = $newName2
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = $newHandoffDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $md1Url, "", "", $newName1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $md2Url, "", "", $newName2)
Style-AsHyperlink $wsOverview.Range("A2")
Style-AsHyperlink $wsOverview.Range("A3")

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newName1
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("D2").Value = "$xlfBase.zh-cn.xlf"
$wsZh.Range("E2").Value = $newHandoffDatetime
$wsZh.Range("F2").Clear()
$wsZh.Range("G2").Clear()
$wsZh.Range("H2").Value = $newHandbackDatetime

$wsZh.Range("A3").Value = $newName2
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("D3").Value = "$xlfBase.zh-cn.xlf"
$wsZh.Range("E3").Value = $newHandoffDatetime
$wsZh.Range("F3").Clear()
$wsZh.Range("G3").Clear()
$wsZh.Range("H3").Value = $newHandbackDatetime

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $md1Url, "", "", $newName1)
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $md1Url, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlf1Url, "", "", "$xlfBase.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $md2Url, "", "", $newName2)
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $md2Url, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhXlf2Url, "", "", "$xlfBase.zh-cn.xlf")
Style-AsHyperlink $wsZh.Range("A2")
Style-AsHyperlink $wsZh.Range("B2")
Style-AsHyperlink $wsZh.Range("D2")
Style-AsHyperlink $wsZh.Range("A3")
Style-AsHyperlink $wsZh.Range("B3")
Style-AsHyperlink $wsZh.Range("D3")

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newName1
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("D2").Value = "$xlfBase.de-de.xlf"
$wsDe.Range("E2").Value = $newDeHandoffDatetime
$wsDe.Range("F2").Clear()
$wsDe.Range("G2").Clear()
$wsDe.Range("H2").Value = $newHandbackDatetime

$wsDe.Range("A3").Value = $newName2
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("D3").Value = "$xlfBase.de-de.xlf"
$wsDe.Range("E3").Value = $newDeHandoffDatetime
$wsDe.Range("F3").Clear()
$wsDe.Range("G3").Clear()
$wsDe.Range("H3").Value = $newHandbackDatetime

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $md1Url, "", "", $newName1)
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $md1Url, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlf1Url, "", "", "$xlfBase.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $md2Url, "", "", $newName2)
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $md2Url, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deXlf2Url, "", "", "$xlfBase.de-de.xlf")
Style-AsHyperlink $wsDe.Range("A2")
Style-AsHyperlink $wsDe.Range("B2")
Style-AsHyperlink $wsDe.Range("D2")
Style-AsHyperlink $wsDe.Range("A3")
Style-AsHyperlink $wsDe.Range("B3")
Style-AsHyperlink $wsDe.Range("D3")
